$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18, column C: replace hardcoded admin email with a {RANDOM_EMAIL} placeholder
$ws.Range("C18").Value = "3.Type `"{RANDOM_EMAIL}`" into `"//input[@placeholder='admin@nesto.com']`""

# Row 21, column C: replace hardcoded mobile number with a {RANDOM_MOBILE} placeholder
$ws.Range("C21").Value = "6.Type `"{RANDOM_MOBILE}`" into `"//input[@placeholder='Your Mobile Number']`""
